$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (stored OOXML width = ColumnWidth + 5/6; target stored width is 15.42578125 for both)
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
$ws.Columns.Item(2).ColumnWidth = 14.592447916666666

$valuesA = @(
  -0.07640048516454101,
  -0.076790567601833004,
  -0.040810053262102741,
  -0.032692064606306914,
  -0.029238048175369791,
  -0.013710455081239914,
  -0.0034512633095911127,
  0.0066020341362693635,
  0.0086779643953436292,
  0.010734404993163693,
  0.013735217269509548,
  0.017239576665915468,
  0.020782863248194339,
  0.028807430927344946,
  0.029835097816138578,
  -0.0060327011669669695,
  -0.0040029870291933278,
  -0.016101344342224877,
  -0.012090392390221805,
  -0.0080157017278796872,
  -0.0040054965516462815,
  -0.031888176335073126,
  -0.040490452868800375,
  -0.020097127786536717,
  -0.011976333495809399,
  -0.0094166001980973135,
  -0.0068412470771415101,
  -0.066223691029193787,
  -0.058821283144288472,
  0.0012865018443020126,
  0.0083582671299087252,
  -0.0040008364970152144
)

$valuesB = @(
  0.076332132558434296,
  0.076698229548283869,
  0.040692064513557114,
  0.032238048131720376,
  0.027686788537423013,
  0.013451263186722073,
  0.0033979657403104824,
  -0.0066779644350769551,
  -0.0087344050312956334,
  -0.010735217317703771,
  -0.013739576718552193,
  -0.01728286329584261,
  -0.020807431019406408,
  -0.028835097832843992,
  -0.029868127883124096,
  0.006002987004738003,
  0.0039999999544866327,
  0.016090392348161231,
  0.012015701682100755,
  0.0080054965053530935,
  0.0039999999533018027,
  0.031707485029301097,
  0.040097127562274792,
  0.019999999772267962,
  0.011916600152858337,
  0.0093412470313705143,
  0.0064050064999783274,
  0.065821283037966971,
  0.0587134974966812,
  -0.0013582672419714115,
  -0.0083659499561221651,
  0.0039999999191113744
)

for ($i = 0; $i -lt 32; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $valuesA[$i]
  $ws.Cells.Item($row, 2).Value = $valuesB[$i]
}
